# Weekly fruit/vegetable price update:
# A new weekly observation is inserted as row 9 (shifting the existing
# rows 9-36 down to 10-37), growing the sheet from A1:R36 to A1:R37.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 9; this shifts rows 9..36 down
# to 10..37 and pushes the sheet dimension from R36 to R37.
$ws.Rows.Item(9).Insert()

# Populate the newly inserted row 9 with the new weekly record. All
# "constant" columns (A, B, C, E, F, G, H, N, O, Q, R) are identical to
# the rest of the data set for this market/product.
$ws.Cells.Item(9, 1).Value = 7
$ws.Cells.Item(9, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(9, 3).Value = "Ñuble"
$ws.Cells.Item(9, 4).Value = 44847
$ws.Cells.Item(9, 5).Value = 16
$ws.Cells.Item(9, 6).Value = 100112043
$ws.Cells.Item(9, 7).Value = "Pepino dulce"
$ws.Cells.Item(9, 8).Value = "Cultivar IV Región"
$ws.Cells.Item(9, 9).Value = "Primera"
$ws.Cells.Item(9, 10).Value = 120
$ws.Cells.Item(9, 11).Value = 17000
$ws.Cells.Item(9, 12).Value = 17000
$ws.Cells.Item(9, 13).Value = 17000
$ws.Cells.Item(9, 14).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(9, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(9, 16).Value = 944
$ws.Cells.Item(9, 17).Value = 18
$ws.Cells.Item(9, 18).Value = "Hortaliza"
